$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.Value = "'" + $newValue
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '68.893.61'
Set-TextValue 'E2' '  -0.09%  '

Set-TextValue 'D3' '3.869.02'
Set-TextValue 'E3' '  +3.18%  '

Set-TextValue 'E4' '  -0.13%  '

Set-TextValue 'D5' '603.13'
Set-TextValue 'E5' '  +0.18%  '

Set-TextValue 'D6' '163.17'
Set-TextValue 'E6' '  -2.50%  '

Set-TextValue 'D7' '3.868.75'
Set-TextValue 'E7' '  +3.24%  '

Set-TextValue 'E8' '  -0.04%  '

Set-TextValue 'D9' '0.531'
Set-TextValue 'E9' '  -1.32%  '

Set-TextValue 'D10' '0.168'
Set-TextValue 'E10' '  -0.40%  '

Set-TextValue 'D11' '6.32'
Set-TextValue 'E11' '  -2.49%  '

Set-TextValue 'D12' '0.460'
Set-TextValue 'E12' '  +0.09%  '

Set-TextValue 'D13' '37.00'
Set-TextValue 'E13' '  -2.31%  '

Set-TextValue 'E14' '  -1.72%  '

Set-TextValue 'D15' '4.510.00'
Set-TextValue 'E15' '  +2.96%  '

Set-TextValue 'D16' '3.849.62'
Set-TextValue 'E16' '  +2.69%  '

Set-TextValue 'D17' '69.074.85'
Set-TextValue 'E17' '  +0.10%  '

Set-TextValue 'D18' '7.56'
Set-TextValue 'E18' '  +2.88%  '

Set-TextValue 'B19' 'Uniswap'
Set-TextValue 'C19' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D19' '11.46'
Set-TextValue 'E19' '  +5.71%  '

Set-TextValue 'B20' 'TRON'
Set-TextValue 'C20' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D20' '0.113'
Set-TextValue 'E20' '  -0.13%  '

Set-TextValue 'D21' '17.18'
Set-TextValue 'E21' '  -0.28%  '

Set-TextValue 'D22' '486.70'
Set-TextValue 'E22' '  -1.17%  '

Set-TextValue 'D23' '0.723'
Set-TextValue 'E23' '  -0.39%  '

Set-TextValue 'D24' '0.0000159'
Set-TextValue 'E24' '  +4.07%  '

Set-TextValue 'D25' '84.13'
Set-TextValue 'E25' '  -0.75%  '

Set-TextValue 'E26' '  -1.78%  '

Set-TextValue 'D27' '12.13'
Set-TextValue 'E27' '  -1.65%  '

Set-TextValue 'D28' '10.05'
Set-TextValue 'E28' '  -0.66%  '

Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  +0.01%  '

Set-TextValue 'D30' '2.98'
Set-TextValue 'E30' '  -0.47%  '

Set-TextValue 'D31' '7.96'
Set-TextValue 'E31' '  -0.79%  '

Set-TextValue 'D32' '4.012.33'
Set-TextValue 'E32' '  +2.97%  '

Set-TextValue 'E33' '  -3.54%  '

Set-TextValue 'D34' '32.40'
Set-TextValue 'E34' '  +2.80%  '

Set-TextValue 'D35' '3.814.80'
Set-TextValue 'E35' '  +3.49%  '

Set-TextValue 'E36' '  -1.63%  '

Set-TextValue 'D37' '1.04'
Set-TextValue 'E37' '  +1.50%  '

Set-TextValue 'E38' '  +4.56%  '

Set-TextValue 'D39' '5.92'
Set-TextValue 'E39' '  +0.45%  '

Set-TextValue 'E40' '  -0.13%  '

Set-TextValue 'B41' 'TheGraph'
Set-TextValue 'C41' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D41' '0.319'
Set-TextValue 'E41' '  -1.49%  '

Set-TextValue 'B42' 'Bittensor'
Set-TextValue 'C42' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D42' '443.21'
Set-TextValue 'E42' '  +3.30%  '

Set-TextValue 'D43' '2.99'
Set-TextValue 'E43' '  +1.00%  '

Set-TextValue 'D44' '48.57'
Set-TextValue 'E44' '  -0.02%  '

Set-TextValue 'D45' '1.98'
Set-TextValue 'E45' '  -0.84%  '

Set-TextValue 'B46' 'Cosmos'
Set-TextValue 'C46' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D46' '8.42'
Set-TextValue 'E46' '  -0.76%  '

Set-TextValue 'B47' 'USDe'
Set-TextValue 'C47' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D47' '1.00'
Set-TextValue 'E47' '  +0.00%  '

Set-TextValue 'D48' '27.19'
Set-TextValue 'E48' '  +16.17%  '

Set-TextValue 'D49' '2.839.27'
Set-TextValue 'E49' '  +1.93%  '

Set-TextValue 'D50' '142.70'
Set-TextValue 'E50' '  +0.80%  '

Set-TextValue 'E51' '  +1.24%  '
